$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB (54) values for rows 1..83 (row1 is the header date row)
$bbVals = @(45986,-1,0.5,-0.4,0.2,-0.2,0,-1.6,1.2,-0.5,1,0.3,-0.4,0.2,-0.3,1.2,-0.7,-0.2,-0.3,-0.6,-0.4,0.5,0,-0.1,0.3,0.1,0,0.4,-0.7,0,-0.1,-0.5,0.2,0.1,-0.4,0,0.3,0.4,-0.5,0.2,0.1,-0.1,-0.1,0.3,0.4,-0.1,-0.8,0.2,-0.9,0.3,0.1,-0.1,-1.9,0.5,2.2,-0.9,0.1,0.2,-0.3,0.7,0.2,0.4,-0.1,0.2,-0.6,-0.7,0,0,0.9,0.5,-0.5,0.4,0.2,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002,0.04000000000000002)

for ($i = 0; $i -lt $bbVals.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 54).Value = $bbVals[$i]
}

# New row 83's date cell in column A
$ws.Range("A83").Value = 46934

# Copy formatting (date number format, bold, border, alignment) from the
# analogous existing cells onto the newly created header/date cells.
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)

$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)

$excel.CutCopyMode = $false
